$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A4").Value = "App"
$ws.Range("B4").Value = "Sarah"
$ws.Range("C4").Value = "Test "
$ws.Range("D4").Value = "2025-10-02 00:17:00"
